# Cap nhat ten bai tap
# Slide 1: merge "Bai " + "19. " runs into a single run "Bai 19. "
# Slide 17 (title "Bai tap"): split into "Bai " + "tap 19.1"
# Slide 18 (title "Bai tap (2)"): split into "Bai tap " + "19.2"

$p = $ppt.ActivePresentation

# --- Slide 1: subtitle "Bai 19. Chia cum van ban (2)" ---
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(2)
$tr1 = $sh1.TextFrame.TextRange
$run1First = $tr1.Characters(1, 4)
$run1First.Text = "Bài 19. "
$run1Second = $tr1.Characters(9, 4)
$run1Second.Text = ""

# --- Slide 17: title "Bai tap" -> "Bai " + "tap 19.1" ---
$s17 = $p.Slides.Item(17)
$sh17 = $s17.Shapes.Item(1)
$tr17 = $sh17.TextFrame.TextRange
$tail17 = $tr17.Characters(5, 3)
$tail17.Text = "tập 19.1"

# --- Slide 18: title "Bai tap (2)" -> "Bai tap " + "19.2" ---
$s18 = $p.Slides.Item(18)
$sh18 = $s18.Shapes.Item(1)
$tr18 = $sh18.TextFrame.TextRange
$tail18 = $tr18.Characters(9, 3)
$tail18.Text = "19.2"
